$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''' + '67.505.22'
$ws.Range('E2').Value = '  +0.80%  '

$ws.Range('D3').Value = '''' + '3.867.00'
$ws.Range('E3').Value = '  +0.35%  '

$ws.Range('E4').Value = '  +0.13%  '

$ws.Range('D5').Value = '''' + '465.48'
$ws.Range('E5').Value = '  +9.78%  '

$ws.Range('D6').Value = '''' + '147.50'
$ws.Range('E6').Value = '  +13.39%  '

$ws.Range('D7').Value = '''' + '0.632'
$ws.Range('E7').Value = '  +3.73%  '

$ws.Range('D8').Value = '''' + '0.999'
$ws.Range('E8').Value = '  +0.08%  '

$ws.Range('D9').Value = '''' + '0.749'
$ws.Range('E9').Value = '  +3.50%  '

$ws.Range('D10').Value = '''' + '0.156'
$ws.Range('E10').Value = '  -1.65%  '

$ws.Range('D11').Value = '''' + '0.0000312'
$ws.Range('E11').Value = '  -7.56%  '

$ws.Range('D12').Value = '''' + '44.05'
$ws.Range('E12').Value = '  +8.10%  '

$ws.Range('D13').Value = '''' + '10.44'
$ws.Range('E13').Value = '  +1.49%  '

$ws.Range('D14').Value = '''' + '4.493.70'
$ws.Range('E14').Value = '  +0.71%  '

$ws.Range('D15').Value = '''' + '14.70'
$ws.Range('E15').Value = '  -6.91%  '

$ws.Range('D16').Value = '''' + '3.880.58'
$ws.Range('E16').Value = '  +0.69%  '

$ws.Range('E17').Value = '  +0.01%  '

$ws.Range('D18').Value = '''' + '20.06'
$ws.Range('E18').Value = '  +1.04%  '

$ws.Range('D19').Value = '''' + '1.16'
$ws.Range('E19').Value = '  +7.53%  '

$ws.Range('D20').Value = '''' + '67.672.10'
$ws.Range('E20').Value = '  +0.64%  '

$ws.Range('D21').Value = '''' + '432.82'
$ws.Range('E21').Value = '  +4.60%  '

$ws.Range('D22').Value = '''' + '14.82'
$ws.Range('E22').Value = '  -0.94%  '

$ws.Range('D23').Value = '''' + '3.28'
$ws.Range('E23').Value = '  +8.09%  '

$ws.Range('D24').Value = '''' + '88.63'
$ws.Range('E24').Value = '  +5.25%  '

$ws.Range('E25').Value = '  +10.27%  '

$ws.Range('B26').Value = 'RenderToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D26').Value = '''' + '10.42'
$ws.Range('E26').Value = '  +12.04%  '

$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = '''' + '37.62'
$ws.Range('E27').Value = '  +0.07%  '

$ws.Range('D28').Value = '''' + '10.25'
$ws.Range('E28').Value = '  +3.58%  '

$ws.Range('D29').Value = '''' + '5.51'
$ws.Range('E29').Value = '  +3.83%  '

$ws.Range('D30').Value = '''' + '748.07'
$ws.Range('E30').Value = '  +3.01%  '

$ws.Range('D31').Value = '''' + '13.78'
$ws.Range('E31').Value = '  +4.73%  '

$ws.Range('D32').Value = '''' + '0.134'
$ws.Range('E32').Value = '  +9.74%  '

$ws.Range('D33').Value = '''' + '2.75'
$ws.Range('E33').Value = '  -0.72%  '

$ws.Range('D34').Value = '''' + '43.32'
$ws.Range('E34').Value = '  +11.88%  '

$ws.Range('D35').Value = '''' + '0.163'
$ws.Range('E35').Value = '  +7.09%  '

$ws.Range('D36').Value = '''' + '57.40'
$ws.Range('E36').Value = '  +3.52%  '

$ws.Range('E37').Value = '  +0.14%  '

$ws.Range('D38').Value = '''' + '5.56'
$ws.Range('E38').Value = '  +0.91%  '

$ws.Range('D39').Value = '''' + '0.0480'
$ws.Range('E39').Value = '  +3.85%  '

$ws.Range('E40').Value = '  +10.48%  '

$ws.Range('D41').Value = '''' + '2.92'
$ws.Range('E41').Value = '  +1.39%  '

$ws.Range('D42').Value = '''' + '2.62'
$ws.Range('E42').Value = '  +13.45%  '

$ws.Range('E43').Value = '  +5.54%  '

$ws.Range('D44').Value = '''' + '0.0₃0675'
$ws.Range('E44').Value = '  -8.35%  '

$ws.Range('E45').Value = '  -0.02%  '

$ws.Range('D46').Value = '''' + '3.28'
$ws.Range('E46').Value = '  +5.63%  '

$ws.Range('D47').Value = '''' + '3.43'
$ws.Range('E47').Value = '  +2.35%  '

$ws.Range('D48').Value = '''' + '2.75'
$ws.Range('E48').Value = '  +7.88%  '

$ws.Range('D49').Value = '''' + '2.13'
$ws.Range('E49').Value = '  +3.60%  '

$ws.Range('D50').Value = '''' + '145.11'
$ws.Range('E50').Value = '  +3.89%  '

$ws.Range('E51').Value = '  +3.53%  '
